$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format ("@") on every Price (column D) cell we touch so that
# numeric-looking strings (e.g. "223.36", "0.0490") are preserved exactly
# as text, matching the inlineStr cells already used throughout the sheet,
# instead of being auto-coerced into floating point numbers by Excel.
$dCells = @("D2","D3","D5","D8","D10","D13","D14","D15","D16","D17","D18","D19","D20","D22","D25","D27","D30","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '34.680.32'
$ws.Range("E2").Value = '  +1.69%  '

$ws.Range("D3").Value = '1.790.97'
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '223.36'
$ws.Range("E5").Value = '  -1.52%  '

$ws.Range("E6").Value = '  -0.55%  '

$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").Value = '32.29'
$ws.Range("E8").Value = '  +3.91%  '

$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("D10").Value = '0.0702'
$ws.Range("E10").Value = '  +6.13%  '

$ws.Range("E11").Value = '  +0.77%  '

$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").Value = '1.793.51'
$ws.Range("E13").Value = '  -0.04%  '

$ws.Range("D14").Value = '10.96'
$ws.Range("E14").Value = '  -2.52%  '

$ws.Range("D15").Value = '34.672.38'
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("D16").Value = '0.631'
$ws.Range("E16").Value = '  -0.35%  '

$ws.Range("D17").Value = '4.29'
$ws.Range("E17").Value = '  +1.87%  '

$ws.Range("D18").Value = '68.92'
$ws.Range("E18").Value = '  -0.95%  '

$ws.Range("D19").Value = '253.40'
$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").Value = '0.0₃0799'
$ws.Range("E20").Value = '  +7.50%  '

$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("D22").Value = '10.65'
$ws.Range("E22").Value = '  +2.39%  '

$ws.Range("E23").Value = '  -1.93%  '

$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("D25").Value = '160.60'
$ws.Range("E25").Value = '  +1.58%  '

$ws.Range("E26").Value = '  -1.54%  '

$ws.Range("D27").Value = '7.09'
$ws.Range("E27").Value = '  +1.25%  '

$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("D30").Value = '0.0526'
$ws.Range("E30").Value = '  +1.90%  '

$ws.Range("E31").Value = '  -3.36%  '

$ws.Range("E32").Value = '  -1.32%  '

$ws.Range("E33").Value = '  -1.00%  '

$ws.Range("E34").Value = '  -0.63%  '

$ws.Range("B35").Value = 'Swop.fi'
$ws.Range("C35").Value = 'https://coinranking.com/coin/yrCr2HW2c+swopfi-swop'
$ws.Range("D35").Value = '440.49'
$ws.Range("E35").Value = '  +733.76%  '

$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '1.435.92'
$ws.Range("E36").Value = '  -3.71%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.0191'
$ws.Range("E37").Value = '  +2.28%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.05'
$ws.Range("E38").Value = '  -1.15%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.632'
$ws.Range("E39").Value = '  -0.06%  '

$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '84.57'
$ws.Range("E40").Value = '  +1.35%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.78'
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '0.922'
$ws.Range("E42").Value = '  +1.80%  '

$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D43").Value = '2.32'
$ws.Range("E43").Value = '  -1.24%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '2.12'
$ws.Range("E44").Value = '  +3.09%  '

$ws.Range("E45").Value = '  +4.12%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.06'
$ws.Range("E46").Value = '  -1.24%  '

$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").Value = '0.0490'
$ws.Range("E47").Value = '  -5.19%  '

$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D48").Value = '1.945.31'
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '105.18'
$ws.Range("E49").Value = '  +7.30%  '

$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.19%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '11.91'
$ws.Range("E51").Value = '  +1.17%  '
